# Update scripts with new TPM (transcripts per million) derived NATMI
# ligand-receptor statistics for the Dkk2-Lrp6 pair. Only the numeric
# columns (Ligand/Receptor/Edge expression, specificity and weight
# measures) that depend on the newly recomputed TPM values change; the
# categorical columns (Sending cluster, Ligand/Receptor symbol, Target
# cluster) keep the exact same text values as before.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0.1764303333333333
$ws.Range("H2").Value = 0.529291
$ws.Range("I2").Value = 0.04559680146739255
$ws.Range("J2").Value = 0.04559680146739255
$ws.Range("M2").Value = 13.17295566666667
$ws.Range("N2").Value = 39.518867
$ws.Range("O2").Value = 0.133784132206724
$ws.Range("P2").Value = 0.133784132206724
$ws.Range("Q2").Value = 2.324108959255222
$ws.Range("R2").Value = 20.916980633297
$ws.Range("S2").Value = 0.006100128515717391
$ws.Range("T2").Value = 0.006100128515717392

$ws.Range("G3").Value = 0.1764303333333333
$ws.Range("H3").Value = 0.529291
$ws.Range("I3").Value = 0.04559680146739255
$ws.Range("J3").Value = 0.04559680146739255
$ws.Range("O3").Value = 0.4382627974978752
$ws.Range("P3").Value = 0.4382627974978752
$ws.Range("Q3").Value = 7.613537400677444
$ws.Range("R3").Value = 68.521836606097
$ws.Range("S3").Value = 0.01998338176805468
$ws.Range("T3").Value = 0.01998338176805468

$ws.Range("G4").Value = 0.1764303333333333
$ws.Range("H4").Value = 0.529291
$ws.Range("I4").Value = 0.04559680146739255
$ws.Range("J4").Value = 0.04559680146739255
$ws.Range("M4").Value = 21.06166566666667
$ws.Range("N4").Value = 63.184997
$ws.Range("O4").Value = 0.2139016281041017
$ws.Range("P4").Value = 0.2139016281041017
$ws.Range("Q4").Value = 3.715916694125222
$ws.Range("R4").Value = 33.443250247127
$ws.Range("S4").Value = 0.009753230070214763
$ws.Range("T4").Value = 0.009753230070214763

$ws.Range("G5").Value = 0.1764303333333333
$ws.Range("H5").Value = 0.529291
$ws.Range("I5").Value = 0.04559680146739255
$ws.Range("J5").Value = 0.04559680146739255
$ws.Range("M5").Value = 21.076417
$ws.Range("N5").Value = 63.229251
$ws.Range("O5").Value = 0.214051442191299
$ws.Range("P5").Value = 0.214051442191299
$ws.Range("Q5").Value = 3.718519276782334
$ws.Range("R5").Value = 33.466673491041
$ws.Range("S5").Value = 0.009760061113405716
$ws.Range("T5").Value = 0.009760061113405716

$ws.Range("I6").Value = 0.9386165989824495
$ws.Range("J6").Value = 0.9386165989824495
$ws.Range("M6").Value = 13.17295566666667
$ws.Range("N6").Value = 39.518867
$ws.Range("O6").Value = 0.133784132206724
$ws.Range("P6").Value = 0.133784132206724
$ws.Range("Q6").Value = 47.84211121827889
$ws.Range("R6").Value = 430.57900096451
$ws.Range("S6").Value = 0.1255720071696936
$ws.Range("T6").Value = 0.1255720071696937

$ws.Range("I7").Value = 0.9386165989824495
$ws.Range("J7").Value = 0.9386165989824495
$ws.Range("O7").Value = 0.4382627974978752
$ws.Range("P7").Value = 0.4382627974978752
$ws.Range("S7").Value = 0.4113607364479896
$ws.Range("T7").Value = 0.4113607364479896

$ws.Range("I8").Value = 0.9386165989824495
$ws.Range("J8").Value = 0.9386165989824495
$ws.Range("M8").Value = 21.06166566666667
$ws.Range("N8").Value = 63.184997
$ws.Range("O8").Value = 0.2139016281041017
$ws.Range("P8").Value = 0.2139016281041017
$ws.Range("Q8").Value = 76.4926700403789
$ws.Range("R8").Value = 688.4340303634101
$ws.Range("S8").Value = 0.2007716186878807
$ws.Range("T8").Value = 0.2007716186878807

$ws.Range("I9").Value = 0.9386165989824495
$ws.Range("J9").Value = 0.9386165989824495
$ws.Range("M9").Value = 21.076417
$ws.Range("N9").Value = 63.229251
$ws.Range("O9").Value = 0.214051442191299
$ws.Range("P9").Value = 0.214051442191299
$ws.Range("Q9").Value = 76.54624457200335
$ws.Range("R9").Value = 688.9162011480302
$ws.Range("S9").Value = 0.2009122366768855
$ws.Range("T9").Value = 0.2009122366768855

$ws.Range("G10").Value = 0.061084
$ws.Range("H10").Value = 0.183252
$ws.Range("I10").Value = 0.01578659955015789
$ws.Range("J10").Value = 0.01578659955015789
$ws.Range("M10").Value = 13.17295566666667
$ws.Range("N10").Value = 39.518867
$ws.Range("O10").Value = 0.133784132206724
$ws.Range("P10").Value = 0.133784132206724
$ws.Range("Q10").Value = 0.8046568239426667
$ws.Range("R10").Value = 7.241911415484
$ws.Range("S10").Value = 0.002111996521312933
$ws.Range("T10").Value = 0.002111996521312933

$ws.Range("G11").Value = 0.061084
$ws.Range("H11").Value = 0.183252
$ws.Range("I11").Value = 0.01578659955015789
$ws.Range("J11").Value = 0.01578659955015789
$ws.Range("O11").Value = 0.4382627974978752
$ws.Range("P11").Value = 0.4382627974978752
$ws.Range("Q11").Value = 2.635971433009333
$ws.Range("R11").Value = 23.723742897084
$ws.Range("S11").Value = 0.006918679281830896
$ws.Range("T11").Value = 0.006918679281830896

$ws.Range("G12").Value = 0.061084
$ws.Range("H12").Value = 0.183252
$ws.Range("I12").Value = 0.01578659955015789
$ws.Range("J12").Value = 0.01578659955015789
$ws.Range("M12").Value = 21.06166566666667
$ws.Range("N12").Value = 63.184997
$ws.Range("O12").Value = 0.2139016281041017
$ws.Range("P12").Value = 0.2139016281041017
$ws.Range("Q12").Value = 1.286530785582667
$ws.Range("R12").Value = 11.578777070244
$ws.Range("S12").Value = 0.003376779346006253
$ws.Range("T12").Value = 0.003376779346006253

$ws.Range("G13").Value = 0.061084
$ws.Range("H13").Value = 0.183252
$ws.Range("I13").Value = 0.01578659955015789
$ws.Range("J13").Value = 0.01578659955015789
$ws.Range("M13").Value = 21.076417
$ws.Range("N13").Value = 63.229251
$ws.Range("O13").Value = 0.214051442191299
$ws.Range("P13").Value = 0.214051442191299
$ws.Range("Q13").Value = 1.287431856028
$ws.Range("R13").Value = 11.586886704252
$ws.Range("S13").Value = 0.00337914440100781
$ws.Range("T13").Value = 0.00337914440100781

